$d = $word.ActiveDocument

# Locate the "5. PARTES INTERESSADAS" table (Nome / Papel no projeto / Assinatura).
$table = $null
for ($i = 1; $i -le $d.Tables.Count; $i++) {
    $candidate = $d.Tables.Item($i)
    if ($candidate.Cell(1, 1).Range.Text -like "*PARTES INTERESSADAS*") {
        $table = $candidate
        break
    }
}

# Add a row for Hugo de Paula (Avaliador), leaving the signature cell blank.
$row1 = $table.Rows.Add()
$row1.Cells.Item(1).Range.Text = "Hugo de Paula"
$row1.Cells.Item(2).Range.Text = "Avaliador"

# Add a row for Pedro Alves (Avaliador), leaving the signature cell blank.
$row2 = $table.Rows.Add()
$row2.Cells.Item(1).Range.Text = "Pedro Alves"
$row2.Cells.Item(2).Range.Text = "Avaliador"

Write-Output "PARTES INTERESSADAS table now has $($table.Rows.Count) rows"
